# Update the header date, then refresh each three-digit x one-digit
# multiplication answer cell in the practice table to the new day's
# problem set. Each old value is a unique, literal string in the
# document, so a straightforward Find/Replace (wrap = find whole
# story, replace-all for that single occurrence) is safe and won't
# clobber unrelated cells.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-19 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-20 Thursday", 2) | Out-Null
$d.Content.Find.Execute("653×7=4571", $true, $false, $false, $false, $false, $true, 1, $false, "446×2=892", 2) | Out-Null
$d.Content.Find.Execute("950×4=3800", $true, $false, $false, $false, $false, $true, 1, $false, "780×7=5460", 2) | Out-Null
$d.Content.Find.Execute("473×9=4257", $true, $false, $false, $false, $false, $true, 1, $false, "133×9=1197", 2) | Out-Null
$d.Content.Find.Execute("696×7=4872", $true, $false, $false, $false, $false, $true, 1, $false, "799×7=5593", 2) | Out-Null
$d.Content.Find.Execute("774×5=3870", $true, $false, $false, $false, $false, $true, 1, $false, "641×3=1923", 2) | Out-Null
$d.Content.Find.Execute("515×2=1030", $true, $false, $false, $false, $false, $true, 1, $false, "298×4=1192", 2) | Out-Null
$d.Content.Find.Execute("631×5=3155", $true, $false, $false, $false, $false, $true, 1, $false, "524×7=3668", 2) | Out-Null
$d.Content.Find.Execute("482×7=3374", $true, $false, $false, $false, $false, $true, 1, $false, "364×3=1092", 2) | Out-Null
$d.Content.Find.Execute("465×3=1395", $true, $false, $false, $false, $false, $true, 1, $false, "304×7=2128", 2) | Out-Null
$d.Content.Find.Execute("118×7=826", $true, $false, $false, $false, $false, $true, 1, $false, "592×4=2368", 2) | Out-Null
$d.Content.Find.Execute("219×9=1971", $true, $false, $false, $false, $false, $true, 1, $false, "791×4=3164", 2) | Out-Null
$d.Content.Find.Execute("890×4=3560", $true, $false, $false, $false, $false, $true, 1, $false, "945×5=4725", 2) | Out-Null
$d.Content.Find.Execute("981×8=7848", $true, $false, $false, $false, $false, $true, 1, $false, "348×8=2784", 2) | Out-Null
$d.Content.Find.Execute("188×8=1504", $true, $false, $false, $false, $false, $true, 1, $false, "994×9=8946", 2) | Out-Null
$d.Content.Find.Execute("610×7=4270", $true, $false, $false, $false, $false, $true, 1, $false, "587×7=4109", 2) | Out-Null
$d.Content.Find.Execute("792×6=4752", $true, $false, $false, $false, $false, $true, 1, $false, "949×2=1898", 2) | Out-Null
$d.Content.Find.Execute("651×2=1302", $true, $false, $false, $false, $false, $true, 1, $false, "855×3=2565", 2) | Out-Null
$d.Content.Find.Execute("733×6=4398", $true, $false, $false, $false, $false, $true, 1, $false, "310×8=2480", 2) | Out-Null
$d.Content.Find.Execute("148×4=592", $true, $false, $false, $false, $false, $true, 1, $false, "825×5=4125", 2) | Out-Null
$d.Content.Find.Execute("114×7=798", $true, $false, $false, $false, $false, $true, 1, $false, "804×6=4824", 2) | Out-Null
$d.Content.Find.Execute("556×4=2224", $true, $false, $false, $false, $false, $true, 1, $false, "892×8=7136", 2) | Out-Null
$d.Content.Find.Execute("233×9=2097", $true, $false, $false, $false, $false, $true, 1, $false, "948×8=7584", 2) | Out-Null
$d.Content.Find.Execute("833×5=4165", $true, $false, $false, $false, $false, $true, 1, $false, "996×5=4980", 2) | Out-Null
$d.Content.Find.Execute("753×7=5271", $true, $false, $false, $false, $false, $true, 1, $false, "103×7=721", 2) | Out-Null
$d.Content.Find.Execute("221×7=1547", $true, $false, $false, $false, $false, $true, 1, $false, "945×7=6615", 2) | Out-Null
